$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value parses as a plain number; these must
# stay stored as text (matching the source data format), so force the
# Text number format before writing the value.
$ws.Range("D2").Value = '61.173.14'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '2.374.10'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.21'
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.71'
$ws.Range("E6").Value = '  -2.01%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.522'
$ws.Range("E8").Value = '  -1.66%  '

$ws.Range("D9").Value = '2.374.58'
$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  +1.57%  '

$ws.Range("E11").Value = '  +1.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.32'
$ws.Range("E12").Value = '  +0.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.347'
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.03'
$ws.Range("E14").Value = '  -1.34%  '

$ws.Range("D15").Value = '2.785.65'
$ws.Range("E15").Value = '  -1.24%  '

$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("D17").Value = '61.071.25'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").Value = '2.367.80'
$ws.Range("E18").Value = '  -0.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.91'
$ws.Range("E19").Value = '  +1.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '319.98'
$ws.Range("E21").Value = '  +0.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.69'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.22'
$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.70'
$ws.Range("E25").Value = '  -10.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.45'
$ws.Range("E26").Value = '  +3.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.14'
$ws.Range("E27").Value = '  +0.81%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0889'
$ws.Range("E28").Value = '  -3.75%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.150'
$ws.Range("E29").Value = '  +2.23%  '

$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '502.77'
$ws.Range("E30").Value = '  -4.89%  '

$ws.Range("E31").Value = '  -4.18%  '

$ws.Range("E32").Value = '  -0.86%  '

$ws.Range("E33").Value = '  -3.53%  '

$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("E36").Value = '  +2.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.42'
$ws.Range("E37").Value = '  -1.54%  '

$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.48'
$ws.Range("E39").Value = '  +2.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '146.74'
$ws.Range("E40").Value = '  +5.13%  '

$ws.Range("E41").Value = '  -0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.19'
$ws.Range("E42").Value = '  +1.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '150.60'
$ws.Range("E43").Value = '  +6.88%  '

$ws.Range("E44").Value = '  -0.15%  '

$ws.Range("E45").Value = '  -0.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0520'
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.27'
$ws.Range("E47").Value = '  -4.01%  '

$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0904'
$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0223'
$ws.Range("E50").Value = '  -0.89%  '

$ws.Range("E51").Value = '  +0.28%  '
